$d = $word.ActiveDocument

# Commit: "Changed license in attribution in quality control section"
#
# In the Quality control section, the image-attribution line at the
# bottom of the "Figure 2: Test-driven Development Cycle" paragraph
# ends with a hyperlinked licence name that reads "CC-BY SA 3.0" and
# needs to become "CC-BY-SA 3.0" (missing hyphen added between "BY"
# and "SA"). That text lives inside a single hyperlink run styled
# with rStyle "InternetLink"; we change only the one space character
# so the run's character formatting is left completely untouched
# (a plain Find/Replace across the whole run text would otherwise
# strip rStyle from the run).

$found = $d.Content.Duplicate
$ok = $found.Find.Execute("CC-BY SA 3.0", $true, $false, $false, $false,
                           $false, $true, 1, $false, "", 0)

if ($ok) {
    $space = $d.Range($found.Start + 5, $found.Start + 6)
    if ($space.Text -eq " ") {
        $space.Text = "-"
    }
}
